# Generate Report for Handback
# Updates timestamps/status produced by a fresh handback-status run for the
# c1037b65-eb75-4124-8cfc-2bab0bbbb831 entry (row 3 and its duplicate row 5).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet - "Latest HO Xliff Generate Date" column (G)
$wsOverview.Range("G3").Value = "2016-08-22 20:14:38"
$wsOverview.Range("G5").Value = "2016-08-22 20:14:38"

# zh-cn sheet - Priority (E), Correspond Handoff Datetime (H), Correspond Handback DateTime (K)
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-08-22 20:14:34"
$wsZhCn.Range("H5").Value = "2016-08-22 20:14:34"
$wsZhCn.Range("K3").Value = "2016-08-22 20:14:50"
$wsZhCn.Range("K5").Value = "2016-08-22 20:14:50"

# de-de sheet - Priority (E), Correspond Handoff Datetime (H), Correspond Handback DateTime (K)
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"
$wsDeDe.Range("H3").Value = "2016-08-22 20:14:38"
$wsDeDe.Range("H5").Value = "2016-08-22 20:14:38"
$wsDeDe.Range("K3").Value = "2016-08-22 20:14:57"
$wsDeDe.Range("K5").Value = "2016-08-22 20:14:57"
